$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose text looks like a number (e.g. "1.003") must be forced to
# Text format before assignment, otherwise Excel auto-converts them to a
# floating point number instead of keeping the literal display string.

$ws.Range('D2').Value = '27.067.85'
$ws.Range('E2').Value = '  -0.44%  '

$ws.Range('D3').Value = '1.891.59'
$ws.Range('E3').Value = '  -0.74%  '

$ws.Range('D4').NumberFormat = "@"
$ws.Range('D4').Value = '1.003'

$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '306.84'
$ws.Range('E5').Value = '  -0.31%  '

$ws.Range('E6').Value = '  +0.15%  '

$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '0.5223'
$ws.Range('E7').Value = '  -0.53%  '

$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '0.3758'
$ws.Range('E8').Value = '  -0.66%  '

$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.07257'
$ws.Range('E9').Value = '  -0.22%  '

$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '21.06'
$ws.Range('E10').Value = '  -0.87%  '

$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '0.8993'
$ws.Range('E11').Value = '  +0.34%  '

$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '0.08189'
$ws.Range('E12').Value = '  +6.54%  '

$ws.Range('D13').Value = '1.937.65'
$ws.Range('E13').Value = '  +1.65%  '

$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '96.23'
$ws.Range('E14').Value = '  +1.31%  '

$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '5.282'
$ws.Range('E15').Value = '  +0.15%  '

$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '1.003'
$ws.Range('E16').Value = '  +0.10%  '

$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '0.000008566'
$ws.Range('E17').Value = '  -0.62%  '

$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '14.57'
$ws.Range('E18').Value = '  +0.59%  '

$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '1.003'
$ws.Range('E19').Value = '  +0.19%  '

$ws.Range('D20').Value = '27.104.04'
$ws.Range('E20').Value = '  -0.55%  '

$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '5.084'
$ws.Range('E21').Value = '  +0.30%  '

$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '10.68'
$ws.Range('E22').Value = '  +0.53%  '

$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '6.401'
$ws.Range('E23').Value = '  -0.53%  '

$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '2.286'
$ws.Range('E25').Value = '  -1.00%  '

$ws.Range('B26').Value = 'EthereumClassic'
$ws.Range('C26').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '18.14'
$ws.Range('E26').Value = '  -0.03%  '

$ws.Range('B27').Value = 'Toncoin'
$ws.Range('C27').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '1.733'
$ws.Range('E27').Value = '  -0.15%  '

$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '114.87'
$ws.Range('E28').Value = '  +0.08%  '

$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '4.782'
$ws.Range('E29').Value = '  -0.38%  '

$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '4.832'
$ws.Range('E30').Value = '  -2.73%  '

$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '0.09223'
$ws.Range('E31').Value = '  -0.03%  '

$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '0.05030'
$ws.Range('E32').Value = '  -0.48%  '

$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '0.7863'
$ws.Range('E33').Value = '  -3.01%  '

$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '1.211'
$ws.Range('E34').Value = '  -2.23%  '

$ws.Range('B35').Value = 'HuobiToken'
$ws.Range('C35').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '2.978'
$ws.Range('E35').Value = '  -0.26%  '

$ws.Range('B36').Value = 'MXToken'
$ws.Range('C36').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '3.417'
$ws.Range('E36').Value = '  +3.27%  '

$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '2.595'
$ws.Range('E37').Value = '  +0.32%  '

$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '0.5728'
$ws.Range('E38').Value = '  +0.83%  '

$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '0.01981'
$ws.Range('E39').Value = '  -0.28%  '

$ws.Range('E40').Value = '  -0.01%  '

$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '9.021'
$ws.Range('E41').Value = '  +0.45%  '

$ws.Range('E42').Value = '  -1.06%  '

$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '116.39'
$ws.Range('E43').Value = '  -2.38%  '

$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '0.1512'
$ws.Range('E44').Value = '  -0.14%  '

$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '0.4849'
$ws.Range('E45').Value = '  +0.41%  '

$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '1.003'
$ws.Range('E46').Value = '  +0.19%  '

$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '10.08'
$ws.Range('E47').Value = '  -1.51%  '

$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '1.619'
$ws.Range('E48').Value = '  -0.09%  '

$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '38.07'
$ws.Range('E49').Value = '  +1.35%  '

$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '63.50'
$ws.Range('E50').Value = '  -0.32%  '

$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '0.05932'
